# Add a "correct answers counter" block (rows 13-15) under the existing
# question table (rows 1-11) on each of the four sheets (ChatGPT, Bard,
# Mistral, Llama), summarising how many ZERO-SHOT / ZERO-SHOT COT answers
# matched the CORRECT ANSWER column, plus the total number of questions.

$wb = $excel.ActiveWorkbook

$labelZeroShot    = "Preguntas acertadas (ZERO-SHOT)"
$labelZeroShotCot = "Preguntas acertadas (ZERO-SHOT COT)"
$labelTotal       = "Preguntas totales"

# Per-sheet counts: [ ZERO-SHOT correct, ZERO-SHOT COT correct, total questions ]
$counts = @{
    "Mistral" = @(2, 2, 9)
    "Llama"   = @(6, 4, 9)
    "Bard"    = @(8, 6, 9)
    "ChatGPT" = @(8, 8, 9)
}

$order = @("Mistral", "Llama", "Bard", "ChatGPT")

foreach ($name in $order) {
    $ws = $wb.Worksheets.Item($name)
    $vals = $counts[$name]

    $ws.Range("A13").Value = $labelZeroShot
    $ws.Range("B13").Value = $vals[0]

    $ws.Range("A14").Value = $labelZeroShotCot
    $ws.Range("B14").Value = $vals[1]

    $ws.Range("A15").Value = $labelTotal
    $ws.Range("B15").Value = $vals[2]

    if ($name -eq "ChatGPT") {
        [void]$ws.Range("A11").Select()
    } elseif ($name -eq "Bard") {
        [void]$ws.Range("A13:B15").Select()
    } else {
        [void]$ws.Range("B15").Select()
    }
}

Write-Output "Added correct-answer counters to ChatGPT, Bard, Mistral and Llama sheets"
